$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 319, pushing existing rows 319:340 down to 320:341
$ws.Rows.Item(319).Insert()

# Populate the new row 319 with this week's record (same market/region/product
# boilerplate as every other row in this sheet)
$ws.Range("A319").Value = 4
$ws.Range("B319").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C319").Value = "Los Lagos"
$ws.Range("D319").Value = 45008
$ws.Range("E319").Value = 10
$ws.Range("F319").Value = "Fruta"
$ws.Range("G319").Value = 100109
$ws.Range("H319").Value = "Uva"
$ws.Range("I319").Value = 100109001
$ws.Range("J319").Value = "Uva"
$ws.Range("K319").Value = "Red Globe"
$ws.Range("L319").Value = "Primera"
$ws.Range("M319").Value = 200
$ws.Range("N319").Value = 14000
$ws.Range("O319").Value = 15000
$ws.Range("P319").Value = 14500
$ws.Range("Q319").Value = "$/caja 18 kilos"
$ws.Range("R319").Value = "Región de O'Higgins"
$ws.Range("S319").Value = 806
$ws.Range("T319").Value = 18
